$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Frontend Bugs")

$rows = @{
    38 = @{ "A"="C-006"; "B"="Notification drop down on customer dashboard"; "C"="Chrome"; "D"="Responsiveness"; "E"="The minimum screen is 280 and the notification drop down goes out of sight"; "F"="Make it visible and responsive"; "G"="must be visible within the screen size"; "H"="irresponsive"; "I"="medium"; "J"="to do"; "K"="bugs\screenshots\customer-page-notification-dropdorn.PNG"; "L"="22nd Feb,2026"; "M"="23rd Feb, 2026"; "O"="Muhammad Noman" }
    39 = @{ "A"="C-007"; "B"="Customer Dashboard cards "; "C"="Chrome"; "D"="Responsiveness"; "E"="The minimum screen is 280 and the cards icon goes out of cards."; "F"="Make it with within the cards"; "G"="Must be adjusted according to the cards size."; "H"="irresponsive"; "I"="medium"; "J"="to do"; "K"="bugs\screenshots\customer-Dasbboard-cards.PNG"; "L"="23rd Feb,2026"; "M"="24th Feb, 2026"; "O"="Muhammad Noman" }
    40 = @{ "A"="C-008"; "B"="Customer Medicine detail page"; "C"="Chrome"; "D"="Responsiveness"; "E"="The minimum screen is 280 and the buy now and cadd to cart buttons and other things goes out of the card"; "F"="Make it responsive"; "G"="Must be adjusted according to the card size"; "H"="irresponsive"; "I"="medium"; "J"="to do"; "K"="bugs\screenshots\customer-medicine-details-page.PNG"; "L"="24th Feb,2026"; "M"="25th Feb, 2026"; "O"="Muhammad Noman" }
    41 = @{ "A"="C-009"; "B"="Customer order detail page"; "C"="Chrome"; "D"="Responsiveness"; "E"="the minimum screen is 280 and the card size miss matched"; "F"="make it responsive"; "G"="Must be adjusted dso that in small screen all cards have same width"; "H"="irresponsive"; "I"="medium"; "J"="to do"; "K"="bugs\screenshots\Customer-order-detail-page.PNG"; "L"="23th Feb,2026"; "M"="24th Feb, 2026"; "O"="Muhammad Noman" }
    42 = @{ "A"="C-0010"; "B"="Customer Appointment page botton size"; "C"="Chrome"; "D"="Responsiveness"; "E"="The buttons on the appointment cards are of unequal lengths make their width equal when the are in small screen"; "F"="make them of equal width"; "G"="Must be adjusted so that in small screen they seem to be of equal width"; "H"="irresponsive"; "I"="medium"; "J"="to do"; "K"="bugs\screenshots\customer-appointment-page-button size.PNG"; "L"="23th Feb,2026"; "M"="24th Feb, 2026"; "O"="Muhammad Noman" }
    43 = @{ "A"="C-0011"; "B"="Customer book Appointment cards and timeline"; "C"="Chrome"; "D"="Responsiveness"; "E"="the min screen size is 280 and the cards of search bar and appointment width mismatched and the time line also goes out of sight"; "F"="Make it visible and responsive"; "G"="must be adjusted to equal widths also time line should be visible"; "H"="irresponsive"; "I"="medium"; "J"="to do"; "K"="bugs\screenshots\Customer-book-appointment responsiveness.PNG"; "L"="22nd Feb, 2026"; "M"="23rd Feb, 2026"; "O"="Muhammad Noman" }
    44 = @{ "A"="C-0012"; "B"="Customer Profile page security tab"; "C"="Chrome"; "D"="Responsiveness"; "E"="the min screen size is 280 and the button goes out of the card"; "F"="make it responsive"; "G"="when the screen size is small the button should be inside the card"; "H"="irresponsive"; "I"="medium"; "J"="to do"; "K"="bugs\screenshots\customer-profile-security-tab.PNG"; "L"="22nd Feb, 2026"; "M"="23rd Feb, 2026"; "O"="Muhammad Noman" }
    45 = @{ "A"="D-001"; "B"="Doctor Registration page logo"; "C"="Chrome"; "D"="Ux"; "E"="the logo of vite is shown instead of philbox logo"; "F"="adjust the logo of philbox"; "G"="must show the logo of philbox"; "H"="vite icon"; "I"="medium"; "J"="to do"; "K"="bugs\screenshots\doctor-registration-page.PNG"; "L"="22nd Feb, 2026"; "M"="23rd Feb, 2026"; "O"="Muhammad Noman" }
    46 = @{ "A"="D-002"; "B"="Doctor login page logo"; "C"="Chrome"; "D"="Ux"; "E"="the logo of vite is shown instead of philbox logo"; "F"="adjust the logo of philbox"; "G"="must show the logo of philbox"; "H"="vite icon"; "I"="medium"; "J"="to do"; "K"="bugs\screenshots\doctor-registration-page.PNG"; "L"="22nd Feb, 2026"; "M"="23rd Feb, 2026"; "O"="Muhammad Noman" }
    47 = @{ "A"="D-003"; "B"="Doctor complete profile"; "C"="Chrome"; "D"="ux"; "E"="the logo of vite is shown instead of philbox logo"; "F"="adjust the logo of philbox"; "G"="must show the logo of philbox"; "H"="vite icon"; "I"="medium"; "J"="to do"; "K"="bugs\screenshots\doctor-registration-page.PNG"; "L"="22nd Feb, 2026"; "M"="23rd Feb, 2026"; "O"="Muhammad Noman" }
    48 = @{ "A"="D-004"; "B"="Doctor Registration page logo"; "C"="Chrome"; "D"="Responsiveness"; "E"="password and confirm password text bxes size mismatched"; "F"="must be adjusted"; "G"="the size of the textboxes must be responsive and well managed"; "H"="irresponsive"; "I"="medium"; "J"="to do"; "K"="bugs\screenshots\doctor-registration-page.PNG"; "L"="22nd Feb, 2026"; "M"="23rd Feb, 2026"; "O"="Muhammad Noman" }
}

foreach ($r in 38..48) {
    $data = $rows[$r]
    foreach ($col in @("A","B","C","D","E","F","G","H","I","J","L","M","O")) {
        if ($data.ContainsKey($col)) {
            $ws.Range("$col$r").Value = $data[$col]
        }
    }
    $ws.Hyperlinks.Add($ws.Range("K$r"), $data["K"])
    $ws.Range("K$r").Style = $ws.Range("K2").Style
}

$ws.Activate()
$ws.Range("O35:O48").Select()
